# To Do List.docx edit:
#   "Job Group A" -> "Ben"
#   "Job Group B" -> "Jacob"   (the _GoBack bookmark moves to sit right
#                               after the new "Jacob" run, before the
#                               following ":" run)
#   "Job Group C" -> "Matthew"
#
# Find/Replace in this engine merges a freshly-edited run into an
# immediately adjacent run that has identical formatting (e.g. "Ben" +
# ":" collapse into a single "Ben:" run). The source diff keeps the
# name and the following ":" as two separate runs, so for the plain
# renames we drop a throwaway bookmark right at the boundary between
# the two runs before doing the replace - that boundary keeps the runs
# from merging - and then remove the marker again afterwards.

$d = $word.ActiveDocument

function Replace-NameKeepingRunSplit($oldText, $newText) {
    $full = $d.Content.Text
    $idx = $full.IndexOf($oldText)
    $endPos = $idx + $oldText.Length
    $boundary = $d.Range($endPos, $endPos)
    $d.Bookmarks.Add("TempRunBoundary", $boundary)

    $find = $d.Content.Find
    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

    $marker = $d.Bookmarks("TempRunBoundary")
    $marker.Delete()
}

Replace-NameKeepingRunSplit "Job Group A" "Ben"
Replace-NameKeepingRunSplit "Job Group B" "Jacob"
Replace-NameKeepingRunSplit "Job Group C" "Matthew"

# Move the "_GoBack" bookmark (it currently sits at the end of the
# "Implement search feature for games." paragraph) so it instead sits
# right after the new "Jacob" run, ahead of the ":" run that follows.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$full = $d.Content.Text
$jacobIdx = $full.IndexOf("Jacob")
$jacobEnd = $jacobIdx + "Jacob".Length
$newBookmarkRange = $d.Range($jacobEnd, $jacobEnd)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
